$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D ("D" and "E"); existing D:K shift to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# New D:E columns inherit default style from the insert; copy number formats
# from column F (the old column D) so D:E match the rest of the quarter columns
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the refreshed quarterly figures (columns D:M) for each line item
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D7:M7").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 290100
$arr[0,1] = 272200
$arr[0,2] = 522000
$arr[0,3] = 279700
$arr[0,4] = 269400
$arr[0,5] = 202700
$arr[0,6] = 370800
$arr[0,7] = 204100
$arr[0,8] = 214400
$arr[0,9] = 203900
$ws.Range("D8:M8").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 216100
$arr[0,1] = 189900
$arr[0,2] = 381600
$arr[0,3] = 202300
$arr[0,4] = 195900
$arr[0,5] = 144500
$arr[0,6] = 268200
$arr[0,7] = 146000
$arr[0,8] = 157300
$arr[0,9] = 134300
$ws.Range("D9:M9").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 74000
$arr[0,1] = 82300
$arr[0,2] = 140400
$arr[0,3] = 77400
$arr[0,4] = 73500
$arr[0,5] = 58200
$arr[0,6] = 102600
$arr[0,7] = 58100
$arr[0,8] = 57100
$arr[0,9] = 69600
$ws.Range("D10:M10").Value = $arr

# Row 11 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = "NA"
$arr[0,1] = "NA"
$arr[0,2] = "NA"
$arr[0,3] = "NA"
$arr[0,4] = "NA"
$arr[0,5] = "NA"
$arr[0,6] = "NA"
$arr[0,7] = "NA"
$arr[0,8] = "NA"
$arr[0,9] = "NA"
$ws.Range("D12:M12").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D13:M13").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 900
$arr[0,1] = 2000
$arr[0,2] = 1500
$arr[0,3] = 800
$arr[0,4] = 7600
$arr[0,5] = 2300
$arr[0,6] = 6200
$arr[0,7] = 4600
$arr[0,8] = 0
$arr[0,9] = 1200
$ws.Range("D14:M14").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 6200
$arr[0,1] = 5000
$arr[0,2] = 10700
$arr[0,3] = 6100
$arr[0,4] = 5200
$arr[0,5] = 2600
$arr[0,6] = 4400
$arr[0,7] = 3000
$arr[0,8] = 2700
$arr[0,9] = 2900
$ws.Range("D15:M15").Value = $arr

# Row 16 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 263800
$arr[0,1] = 240700
$arr[0,2] = 487700
$arr[0,3] = 262800
$arr[0,4] = 253400
$arr[0,5] = 195500
$arr[0,6] = 360400
$arr[0,7] = 203800
$arr[0,8] = 212000
$arr[0,9] = 183900
$ws.Range("D17:M17").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 26300
$arr[0,1] = 31500
$arr[0,2] = 34300
$arr[0,3] = 16900
$arr[0,4] = 16000
$arr[0,5] = 7200
$arr[0,6] = 10400
$arr[0,7] = 300
$arr[0,8] = 2400
$arr[0,9] = 20000
$ws.Range("D18:M18").Value = $arr

# Row 19 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -600
$arr[0,1] = 300
$arr[0,2] = -100
$arr[0,3] = -1600
$arr[0,4] = 5600
$arr[0,5] = 3400
$arr[0,6] = 2100
$arr[0,7] = -300
$arr[0,8] = -300
$arr[0,9] = 0
$ws.Range("D20:M20").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 39600
$arr[0,1] = 41800
$arr[0,2] = 61100
$arr[0,3] = 29000
$arr[0,4] = 34000
$arr[0,5] = 17300
$arr[0,6] = 31000
$arr[0,7] = 9100
$arr[0,8] = 11000
$arr[0,9] = 29500
$ws.Range("D21:M21").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 3800
$arr[0,1] = 5600
$arr[0,2] = 13300
$arr[0,3] = 6700
$arr[0,4] = 11000
$arr[0,5] = 9900
$arr[0,6] = 12700
$arr[0,7] = 4700
$arr[0,8] = 5100
$arr[0,9] = 4600
$ws.Range("D22:M22").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 21900
$arr[0,1] = 26200
$arr[0,2] = 20900
$arr[0,3] = 8600
$arr[0,4] = 10600
$arr[0,5] = 700
$arr[0,6] = -200
$arr[0,7] = -4700
$arr[0,8] = -3000
$arr[0,9] = 15400
$ws.Range("D23:M23").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 7300
$arr[0,1] = 4200
$arr[0,2] = 5500
$arr[0,3] = 2300
$arr[0,4] = 5200
$arr[0,5] = 700
$arr[0,6] = -100
$arr[0,7] = -1800
$arr[0,8] = 900
$arr[0,9] = 1800
$ws.Range("D24:M24").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D25:M25").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 14600
$arr[0,1] = 22000
$arr[0,2] = 15400
$arr[0,3] = 6300
$arr[0,4] = 5400
$arr[0,5] = 0
$arr[0,6] = -100
$arr[0,7] = -3000
$arr[0,8] = -3900
$arr[0,9] = 13700
$ws.Range("D26:M26").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 14400
$arr[0,1] = 21500
$arr[0,2] = 14100
$arr[0,3] = 5800
$arr[0,4] = 5000
$arr[0,5] = -600
$arr[0,6] = -600
$arr[0,7] = -2900
$arr[0,8] = -3300
$arr[0,9] = 15000
$ws.Range("D27:M27").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D28:M28").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 33300
$arr[0,1] = 700
$arr[0,2] = 4000
$arr[0,3] = "NA"
$arr[0,4] = 21600
$arr[0,5] = 2100
$arr[0,6] = 500
$arr[0,7] = "NA"
$arr[0,8] = "NA"
$arr[0,9] = "NA"
$ws.Range("D29:M29").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D30:M30").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D31:M31").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 600
$arr[0,1] = -300
$arr[0,2] = 100
$arr[0,3] = 1600
$arr[0,4] = -5600
$arr[0,5] = -3400
$arr[0,6] = -2100
$arr[0,7] = 300
$arr[0,8] = 300
$arr[0,9] = 0
$ws.Range("D32:M32").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 47700
$arr[0,1] = 22200
$arr[0,2] = 18100
$arr[0,3] = 5800
$arr[0,4] = 26600
$arr[0,5] = 1500
$arr[0,6] = -100
$arr[0,7] = -2900
$arr[0,8] = -3300
$arr[0,9] = 15000
$ws.Range("D33:M33").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D34:M34").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 47700
$arr[0,1] = 22200
$arr[0,2] = 18100
$arr[0,3] = 5800
$arr[0,4] = 26600
$arr[0,5] = 1500
$arr[0,6] = -100
$arr[0,7] = -2900
$arr[0,8] = -3300
$arr[0,9] = 15000
$ws.Range("D35:M35").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D38:M38").Value = $arr

# Row 39 stays blank in D:M
# Row 40 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 118100
$arr[0,1] = 157200
$arr[0,2] = 148500
$arr[0,3] = 142800
$arr[0,4] = 122600
$arr[0,5] = 124700
$arr[0,6] = 234400
$arr[0,7] = 243800
$arr[0,8] = 282000
$arr[0,9] = 267300
$ws.Range("D41:M41").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D42:M42").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 249300
$arr[0,1] = 211500
$arr[0,2] = 245600
$arr[0,3] = 245000
$arr[0,4] = 232900
$arr[0,5] = 237200
$arr[0,6] = 183300
$arr[0,7] = 165300
$arr[0,8] = 169500
$arr[0,9] = 152200
$ws.Range("D43:M43").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 233100
$arr[0,1] = 188200
$arr[0,2] = 233500
$arr[0,3] = 223100
$arr[0,4] = 173700
$arr[0,5] = 213600
$arr[0,6] = 184300
$arr[0,7] = 185300
$arr[0,8] = 169700
$arr[0,9] = 182200
$ws.Range("D44:M44").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 61200
$arr[0,1] = 110900
$arr[0,2] = 33400
$arr[0,3] = 39100
$arr[0,4] = 104800
$arr[0,5] = 46800
$arr[0,6] = 33700
$arr[0,7] = 38700
$arr[0,8] = 31800
$arr[0,9] = 26100
$ws.Range("D45:M45").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 661700
$arr[0,1] = 667800
$arr[0,2] = 661000
$arr[0,3] = 650000
$arr[0,4] = 634000
$arr[0,5] = 622200
$arr[0,6] = 635700
$arr[0,7] = 633000
$arr[0,8] = 653000
$arr[0,9] = 627700
$ws.Range("D46:M46").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D47:M47").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 361100
$arr[0,1] = 288400
$arr[0,2] = 304000
$arr[0,3] = 305000
$arr[0,4] = 285000
$arr[0,5] = 293100
$arr[0,6] = 259800
$arr[0,7] = 254600
$arr[0,8] = 251000
$arr[0,9] = 258100
$ws.Range("D48:M48").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 851100
$arr[0,1] = 727500
$arr[0,2] = 763000
$arr[0,3] = 775000
$arr[0,4] = 746100
$arr[0,5] = 756400
$arr[0,6] = 325300
$arr[0,7] = 327400
$arr[0,8] = 311400
$arr[0,9] = 316000
$ws.Range("D49:M49").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D50:M50").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D51:M51").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 23800
$arr[0,1] = 54800
$arr[0,2] = 21300
$arr[0,3] = 22100
$arr[0,4] = 59600
$arr[0,5] = 21300
$arr[0,6] = 18500
$arr[0,7] = 17600
$arr[0,8] = 17600
$arr[0,9] = 16000
$ws.Range("D52:M52").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D53:M53").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1897700
$arr[0,1] = 1738500
$arr[0,2] = 1749300
$arr[0,3] = 1752100
$arr[0,4] = 1724700
$arr[0,5] = 1693000
$arr[0,6] = 1239300
$arr[0,7] = 1232600
$arr[0,8] = 1233100
$arr[0,9] = 1217900
$ws.Range("D54:M54").Value = $arr

# Row 55 stays blank in D:M
# Row 56 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 125500
$arr[0,1] = 94800
$arr[0,2] = 128600
$arr[0,3] = 124400
$arr[0,4] = 105400
$arr[0,5] = 109900
$arr[0,6] = 82400
$arr[0,7] = 88900
$arr[0,8] = 80000
$arr[0,9] = 65700
$ws.Range("D57:M57").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 11200
$arr[0,1] = 9700
$arr[0,2] = 66800
$arr[0,3] = 63000
$arr[0,4] = 58900
$arr[0,5] = 244300
$arr[0,6] = 6600
$arr[0,7] = 5100
$arr[0,8] = 6500
$arr[0,9] = 6700
$ws.Range("D58:M58").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 229900
$arr[0,1] = 208300
$arr[0,2] = 197600
$arr[0,3] = 194600
$arr[0,4] = 223300
$arr[0,5] = 196200
$arr[0,6] = 151900
$arr[0,7] = 158200
$arr[0,8] = 175100
$arr[0,9] = 152100
$ws.Range("D59:M59").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 366600
$arr[0,1] = 312800
$arr[0,2] = 393000
$arr[0,3] = 382000
$arr[0,4] = 387600
$arr[0,5] = 550500
$arr[0,6] = 240900
$arr[0,7] = 252100
$arr[0,8] = 261500
$arr[0,9] = 224600
$ws.Range("D60:M60").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 533200
$arr[0,1] = 500500
$arr[0,2] = 443900
$arr[0,3] = 449200
$arr[0,4] = 439200
$arr[0,5] = 304000
$arr[0,6] = 237900
$arr[0,7] = 237200
$arr[0,8] = 233700
$arr[0,9] = 232100
$ws.Range("D61:M61").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 108900
$arr[0,1] = 84200
$arr[0,2] = 92300
$arr[0,3] = 92300
$arr[0,4] = 92700
$arr[0,5] = 106100
$arr[0,6] = 40100
$arr[0,7] = 40400
$arr[0,8] = 39200
$arr[0,9] = 48900
$ws.Range("D62:M62").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D63:M63").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D64:M64").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D65:M65").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1013200
$arr[0,1] = 901700
$arr[0,2] = 933000
$arr[0,3] = 927100
$arr[0,4] = 922500
$arr[0,5] = 963200
$arr[0,6] = 520800
$arr[0,7] = 531100
$arr[0,8] = 535800
$arr[0,9] = 507600
$ws.Range("D66:M66").Value = $arr

# Row 67 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D68:M68").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D69:M69").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D70:M70").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D71:M71").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 453900
$arr[0,1] = 406900
$arr[0,2] = 384700
$arr[0,3] = 372400
$arr[0,4] = 364300
$arr[0,5] = 337700
$arr[0,6] = 336200
$arr[0,7] = 333400
$arr[0,8] = 336300
$arr[0,9] = 339600
$ws.Range("D72:M72").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D73:M73").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D74:M74").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D75:M75").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 884500
$arr[0,1] = 836800
$arr[0,2] = 816300
$arr[0,3] = 825000
$arr[0,4] = 802200
$arr[0,5] = 729800
$arr[0,6] = 718500
$arr[0,7] = 701500
$arr[0,8] = 697300
$arr[0,9] = 710300
$ws.Range("D76:M76").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D77:M77").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D80:M80").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 47700
$arr[0,1] = 22200
$arr[0,2] = 18100
$arr[0,3] = 5800
$arr[0,4] = 26600
$arr[0,5] = 1500
$arr[0,6] = -100
$arr[0,7] = -2900
$arr[0,8] = -3300
$arr[0,9] = 15000
$ws.Range("D81:M81").Value = $arr

# Row 82 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 13900
$arr[0,1] = 10000
$arr[0,2] = 26900
$arr[0,3] = 13700
$arr[0,4] = 13200
$arr[0,5] = 10200
$arr[0,6] = 18500
$arr[0,7] = 9200
$arr[0,8] = 8900
$arr[0,9] = 9400
$ws.Range("D83:M83").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D84:M84").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D85:M85").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D86:M86").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D87:M87").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D88:M88").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 23700
$arr[0,1] = 18400
$arr[0,2] = 46700
$arr[0,3] = 23000
$arr[0,4] = 29500
$arr[0,5] = 20400
$arr[0,6] = -2900
$arr[0,7] = 1200
$arr[0,8] = 24200
$arr[0,9] = 59800
$ws.Range("D89:M89").Value = $arr

# Row 90 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -9200
$arr[0,1] = -7300
$arr[0,2] = -19100
$arr[0,3] = -6600
$arr[0,4] = -11200
$arr[0,5] = -5000
$arr[0,6] = -16800
$arr[0,7] = -8400
$arr[0,8] = -4400
$arr[0,9] = -4000
$ws.Range("D91:M91").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D92:M92").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D93:M93").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -89000
$arr[0,1] = -8000
$arr[0,2] = -30900
$arr[0,3] = -19000
$arr[0,4] = -11900
$arr[0,5] = -429200
$arr[0,6] = -38900
$arr[0,7] = -31400
$arr[0,8] = -4300
$arr[0,9] = -3600
$ws.Range("D94:M94").Value = $arr

# Row 95 stays blank in D:M
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D96:M96").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D97:M97").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D98:M98").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D99:M99").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 30900
$arr[0,1] = 1200
$arr[0,2] = 6100
$arr[0,3] = 10000
$arr[0,4] = -21000
$arr[0,5] = 300100
$arr[0,6] = -3900
$arr[0,7] = -3200
$arr[0,8] = -1100
$arr[0,9] = -2300
$ws.Range("D100:M100").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -4700
$arr[0,1] = -2900
$arr[0,2] = -3800
$arr[0,3] = 3900
$arr[0,4] = 2300
$arr[0,5] = 1300
$arr[0,6] = 3600
$arr[0,7] = 700
$arr[0,8] = -4100
$arr[0,9] = 200
$ws.Range("D101:M101").Value = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -39100
$arr[0,1] = 8700
$arr[0,2] = 18100
$arr[0,3] = 17900
$arr[0,4] = -1100
$arr[0,5] = -107400
$arr[0,6] = -42100
$arr[0,7] = -32700
$arr[0,8] = 14700
$arr[0,9] = 54100
$ws.Range("D102:M102").Value = $arr

